# Apply updated leve-profit values (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1522.125
$ws.Range("I38").Value = 548.7273
$ws.Range("J38").Value = 3663.6
$ws.Range("K38").Value = 1646.1819
$ws.Range("L38").Value = 10990.8
$ws.Range("M38").Value = -1274.1819
$ws.Range("N38").Value = -11734.8
$ws.Range("H58").Value = 1364.4445
$ws.Range("I58").Value = 182.85715
$ws.Range("J58").Value = 5500
$ws.Range("K58").Value = 548.5714499999999
$ws.Range("L58").Value = 16500
$ws.Range("M58").Value = -398.5714499999999
$ws.Range("N58").Value = -16800
$ws.Range("H80").Value = 489.5
$ws.Range("I80").Value = 463
$ws.Range("J80").Value = 542.5
$ws.Range("K80").Value = 1389
$ws.Range("L80").Value = 1627.5
$ws.Range("M80").Value = -391
$ws.Range("N80").Value = -3623.5
$ws.Range("H83").Value = 489.5
$ws.Range("I83").Value = 463
$ws.Range("J83").Value = 542.5
$ws.Range("K83").Value = 4167
$ws.Range("L83").Value = 4882.5
$ws.Range("M83").Value = 825
$ws.Range("N83").Value = -14866.5
$ws.Range("H88").Value = 12977.556
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 16114
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 16114
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -16926
$ws.Range("H91").Value = 12977.556
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 16114
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 16114
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -18922
$ws.Range("H92").Value = 121213070
$ws.Range("I92").Value = 920.5714
$ws.Range("J92").Value = 333334340
$ws.Range("K92").Value = 920.5714
$ws.Range("L92").Value = 333334340
$ws.Range("M92").Value = 327.4286
$ws.Range("N92").Value = -333336836
$ws.Range("H129").Value = 1063.9512
$ws.Range("I129").Value = 440
$ws.Range("J129").Value = 1192.4117
$ws.Range("K129").Value = 1320
$ws.Range("L129").Value = 3577.2351
$ws.Range("M129").Value = 3680
$ws.Range("N129").Value = -13577.2351
$ws.Range("H132").Value = 1465.75
$ws.Range("I132").Value = 961.28125
$ws.Range("K132").Value = 2883.84375
$ws.Range("M132").Value = -353.84375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 799.2632
$ws.Range("I2").Value = 567.0625
$ws.Range("J2").Value = 2037.6666
$ws.Range("K2").Value = 567.0625
$ws.Range("L2").Value = 2037.6666
$ws.Range("M2").Value = -454.0625
$ws.Range("N2").Value = -2263.6666
$ws.Range("H116").Value = 799.2632
$ws.Range("I116").Value = 567.0625
$ws.Range("J116").Value = 2037.6666
$ws.Range("K116").Value = 567.0625
$ws.Range("L116").Value = 2037.6666
$ws.Range("M116").Value = 1726.9375
$ws.Range("N116").Value = -6625.6666
$ws.Range("H122").Value = 3370.8572
$ws.Range("I122").Value = 3464.4
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 10393.2
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -7943.200000000001
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 799.2632
$ws.Range("I3").Value = 567.0625
$ws.Range("J3").Value = 2037.6666
$ws.Range("K3").Value = 567.0625
$ws.Range("L3").Value = 2037.6666
$ws.Range("M3").Value = -453.0625
$ws.Range("N3").Value = -2265.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1974.6364
$ws.Range("I16").Value = 1665.4286
$ws.Range("K16").Value = 1665.4286
$ws.Range("M16").Value = -1378.4286
$ws.Range("H113").Value = 1974.6364
$ws.Range("I113").Value = 1665.4286
$ws.Range("K113").Value = 1665.4286
$ws.Range("M113").Value = 504.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3025.6667
$ws.Range("J70").Value = 3825.3333
$ws.Range("L70").Value = 11475.9999
$ws.Range("N70").Value = -12105.9999
$ws.Range("H73").Value = 3025.6667
$ws.Range("J73").Value = 3825.3333
$ws.Range("L73").Value = 11475.9999
$ws.Range("N73").Value = -13659.9999
$ws.Range("H87").Value = 5996
$ws.Range("J87").Value = 5500
$ws.Range("L87").Value = 16500
$ws.Range("N87").Value = -18996
$ws.Range("H90").Value = 5996
$ws.Range("J90").Value = 5500
$ws.Range("L90").Value = 49500
$ws.Range("N90").Value = -61980
$ws.Range("H114").Value = 5841.364
$ws.Range("J114").Value = 8967.357
$ws.Range("L114").Value = 26902.071
$ws.Range("N114").Value = -33410.071
$ws.Range("H122").Value = 9775.817999999999
$ws.Range("I122").Value = 550.6667
$ws.Range("J122").Value = 20846
$ws.Range("K122").Value = 4956.0003
$ws.Range("L122").Value = 187614
$ws.Range("M122").Value = -2506.0003
$ws.Range("N122").Value = -192514

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6930.0454
$ws.Range("I80").Value = 9375.357
$ws.Range("J80").Value = 2650.75
$ws.Range("K80").Value = 9375.357
$ws.Range("L80").Value = 2650.75
$ws.Range("M80").Value = -8377.357
$ws.Range("N80").Value = -4646.75
$ws.Range("H83").Value = 6930.0454
$ws.Range("I83").Value = 9375.357
$ws.Range("J83").Value = 2650.75
$ws.Range("K83").Value = 46876.785
$ws.Range("L83").Value = 13253.75
$ws.Range("M83").Value = -41884.785
$ws.Range("N83").Value = -23237.75
$ws.Range("H122").Value = 1159853.5
$ws.Range("I122").Value = 2851039.2
$ws.Range("J122").Value = 2726.3157
$ws.Range("K122").Value = 8553117.600000001
$ws.Range("L122").Value = 8178.9471
$ws.Range("M122").Value = -8550667.600000001
$ws.Range("N122").Value = -13078.9471
$ws.Range("H123").Value = 20903.148
$ws.Range("J123").Value = 21360.96
$ws.Range("L123").Value = 21360.96
$ws.Range("N123").Value = -26260.96
$ws.Range("H132").Value = 2662.3257
$ws.Range("I132").Value = 2938.182
$ws.Range("J132").Value = 2373.3333
$ws.Range("K132").Value = 8814.545999999998
$ws.Range("L132").Value = 7119.999899999999
$ws.Range("M132").Value = -6284.545999999998
$ws.Range("N132").Value = -12179.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1796.8
$ws.Range("I16").Value = 1918.6666
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1918.6666
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -1748.6666
$ws.Range("N16").Value = -1040
$ws.Range("H61").Value = 929.8889
$ws.Range("I61").Value = 961.125
$ws.Range("J61").Value = 884.4545000000001
$ws.Range("K61").Value = 961.125
$ws.Range("L61").Value = 884.4545000000001
$ws.Range("M61").Value = -759.125
$ws.Range("N61").Value = -1288.4545
$ws.Range("H113").Value = 929.8889
$ws.Range("I113").Value = 961.125
$ws.Range("J113").Value = 884.4545000000001
$ws.Range("K113").Value = 961.125
$ws.Range("L113").Value = 884.4545000000001
$ws.Range("M113").Value = 1208.875
$ws.Range("N113").Value = -5224.4545

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1266.9166
$ws.Range("I81").Value = 1266.9166
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2533.8332
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1472.8332
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1266.9166
$ws.Range("I84").Value = 1266.9166
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 12669.166
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -7365.166000000001
$ws.Range("N84").ClearContents()
